# Add quiz results rows for the "added quiz without payment" entries.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("Chinmay Mhatre",          2, "chinmayrmhatre@gmail.com",     "13/12/2024, 4:10:39 pm", "05779eb0"),
    @("Chinmay Mhatre",          2, "chinmayrmhatre@gmail.com",     "13/12/2024, 4:39:25 pm", "8f08a72f"),
    @("Khushil Girish Bhimani",  2, "Khushilbhimani2@gmail.com",    "15/12/2024, 7:12:05 pm", "f8ab2fd2"),
    @("Khushil Bhimani",         2, "khushilbhimani302@gmail.com",  "15/12/2024, 7:17:29 pm", "dc89af5f"),
    @("Khushil Girish Bhimani",  2, "khushilbhimani2@gmail.com",    "15/12/2024, 7:31:07 pm", "24abf7bc"),
    @("Rohit Gautam",            0, "rohit@hacktify.in",            "15/12/2024, 7:42:21 pm", "b609cfb8")
)

$startRow = 24
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
